# Auto update Excel log
# Appends newly-logged sensor readings to the ALERTS, PIR, Humidity,
# Temperature and mmWave sheets of the SeniorConnect master log.

$wb = $excel.ActiveWorkbook

function Add-LogRow($ws, $row, $date, $timestamp, $hour, $location, $value, $status, $valueIsPercent) {
    # Column A always holds a date-looking string (e.g. "2026-01-30").
    # Excel's COM layer auto-parses such strings into date serials, so we
    # force the cell to Text format first to preserve the literal string.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $date

    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location

    if ($valueIsPercent) {
        # Percent-looking strings (e.g. "87.7%") get auto-parsed into
        # numeric percentages, so keep them literal as text too.
        $ws.Cells.Item($row, 5).NumberFormat = "@"
    }
    $ws.Cells.Item($row, 5).Value = $value

    $ws.Cells.Item($row, 6).Value = $status
}

# ---------------------------------------------------------------------
# ALERTS sheet: one new critical fall alert (row 8)
# ---------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $wsAlerts 8 "2026-01-30" "16:56:57" "16:00" "Living Room" "CRITICAL EMERGENCY" "FALL_DETECTED" $false

# ---------------------------------------------------------------------
# PIR sheet: rows 88-100, Bathroom / No Motion / Inactive
# ---------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")
$pirTimes = @(
    "16:56:08","16:56:09","16:56:13","16:56:18","16:56:23","16:56:28",
    "16:56:33","16:56:38","16:56:43","16:56:48","16:56:53","16:56:58","16:57:03"
)
$r = 88
foreach ($t in $pirTimes) {
    Add-LogRow $wsPir $r "2026-01-30" $t "16:00" "Bathroom" "No Motion" "Inactive" $false
    $r++
}

# ---------------------------------------------------------------------
# Humidity sheet: rows 52-62, Bathroom / xx.x% / Active
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @{ Time = "16:56:08"; Value = "87.7%" },
    @{ Time = "16:56:09"; Value = "87.7%" },
    @{ Time = "16:56:13"; Value = "86.3%" },
    @{ Time = "16:56:23"; Value = "87.7%" },
    @{ Time = "16:56:28"; Value = "86.8%" },
    @{ Time = "16:56:33"; Value = "87.7%" },
    @{ Time = "16:56:43"; Value = "87.8%" },
    @{ Time = "16:56:48"; Value = "86.8%" },
    @{ Time = "16:56:53"; Value = "87.7%" },
    @{ Time = "16:56:58"; Value = "87.7%" },
    @{ Time = "16:57:03"; Value = "87.7%" }
)
$r = 52
foreach ($row in $humidityRows) {
    Add-LogRow $wsHumidity $r "2026-01-30" $row.Time "16:00" "Bathroom" $row.Value "Active" $true
    $r++
}

# ---------------------------------------------------------------------
# Temperature sheet: rows 14-24, Bathroom / xx.xC / Active
# ---------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @{ Time = "16:56:09"; Value = "22.6C" },
    @{ Time = "16:56:09"; Value = "22.6C" },
    @{ Time = "16:56:13"; Value = "22.7C" },
    @{ Time = "16:56:23"; Value = "22.6C" },
    @{ Time = "16:56:28"; Value = "22.6C" },
    @{ Time = "16:56:33"; Value = "22.6C" },
    @{ Time = "16:56:43"; Value = "22.7C" },
    @{ Time = "16:56:48"; Value = "22.6C" },
    @{ Time = "16:56:53"; Value = "22.6C" },
    @{ Time = "16:56:58"; Value = "22.6C" },
    @{ Time = "16:57:03"; Value = "22.7C" }
)
$r = 14
foreach ($row in $temperatureRows) {
    Add-LogRow $wsTemperature $r "2026-01-30" $row.Time "16:00" "Bathroom" $row.Value "Active" $false
    $r++
}

# ---------------------------------------------------------------------
# mmWave sheet: rows 22-24, Living Room / PRESENCE_DETECTED / Active
# ---------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")
$mmWaveTimes = @("16:56:10","16:56:19","16:56:29")
$r = 22
foreach ($t in $mmWaveTimes) {
    Add-LogRow $wsMmWave $r "2026-01-30" $t "16:00" "Living Room" "PRESENCE_DETECTED" "Active" $false
    $r++
}
